$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.985.06"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "1.711.69"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3962"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4039"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08842"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.483"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.021"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001358"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "1.716.81"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07196"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.417"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").Value = "24.987.66"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.990"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.356"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.213"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "151.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.437"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.475"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +24.81%  "
$ws.Range("D33").Value = "1.904.92"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08557"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03164"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.64%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.223"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.046"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09527"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8272"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.483"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.674"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7410"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.261"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.381"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08727"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.89%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  +0.18%  "
